$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.128.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.391.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.385.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.544"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.81%  "

$ws.Range("E11").Value = "  +2.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.971.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("E14").Value = "  -3.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000189"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.171.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.388.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.529"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.05%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.96%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.965.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.52%  "

$ws.Range("E41").Value = "  +1.27%  "

$ws.Range("E42").Value = "  -2.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +21.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.832"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.47%  "
